$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete rows (old rows 4 and 5 -> MuSCs/ECs pairings no longer exist
# with the refreshed TPM numbers).
$ws.Rows("4:5").Delete()

# Row 2: FAPs -> Cxcl5/Cxcr1 -> Resolving-Mac (target cluster renamed from ECs), with
# refreshed TPM-derived statistics.
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl5"
$ws.Range("C2").Value = "Cxcr1"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.702873666666666
$ws.Range("H2").Value = 14.108621
$ws.Range("I2").Value = 0.9961712500318616
$ws.Range("J2").Value = 0.9961712500318616
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.009727666666666667
$ws.Range("N2").Value = 0.029183
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.04574798740477778
$ws.Range("R2").Value = 0.411731886643
$ws.Range("S2").Value = 0.9961712500318616
$ws.Range("T2").Value = 0.9961712500318616

# Row 3: sending cluster renamed from FAPs to Resolving-Mac, target cluster also
# Resolving-Mac, with refreshed TPM-derived statistics.
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Cxcl5"
$ws.Range("C3").Value = "Cxcr1"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01807533333333334
$ws.Range("H3").Value = 0.054226
$ws.Range("I3").Value = 0.003828749968138469
$ws.Range("J3").Value = 0.003828749968138468
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.009727666666666667
$ws.Range("N3").Value = 0.029183
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.0001758308175555556
$ws.Range("R3").Value = 0.001582477358
$ws.Range("S3").Value = 0.003828749968138469
$ws.Range("T3").Value = 0.003828749968138468
